$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("D1").Value = "ncc_anh"
$ws.Range("E1").Value = "ncc_ngaytao"
$ws.Range("F1").Value = "ncc_ngaycapnhat"

# New data values in column D
$ws.Range("D2").Value = "/suppliers/anh1"
$ws.Range("D3").Value = "/suppliers/anh2"
$ws.Range("D4").Value = "/suppliers/anh3"
$ws.Range("D5").Value = "/suppliers/anh4"

# Column widths to match bestFit widths from diff (closest achievable
# quantization under this runtime's width model)
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 11
$ws.Columns.Item(6).ColumnWidth = 15.3333333333333

# Update selection to D2 as in diff
$ws.Range("D2").Select()

$wb.Save()
